$wb = $excel.ActiveWorkbook

# --- "About" sheet updates ---
$about = $wb.Worksheets.Item("About")
$about.Range("B1").Value = "Louisiana"
$about.Range("C1").Value = 44463

# --- "BDPbES" sheet updates ---
$bdpbes = $wb.Worksheets.Item("BDPbES")
# Onshore wind priority: 2 -> 1
$bdpbes.Range("B6").Value = 1
# Solar PV priority: 2 -> 1
$bdpbes.Range("B7").Value = 1

# Make BDPbES the active sheet/tab with B8 selected (matches new saved view state)
$bdpbes.Activate()
$bdpbes.Range("B8").Select()
